$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il23a"
$ws.Range("C2").Value = "Il23r"
$ws.Range("D2").Value = "Neutrophils"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5663796666666666
$ws.Range("H2").Value = 1.699139
$ws.Range("I2").Value = 0.005761481588340165
$ws.Range("J2").Value = 0.005761481588340165
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2324893333333333
$ws.Range("N2").Value = 0.697468
$ws.Range("O2").Value = 0.920195895276243
$ws.Range("P2").Value = 0.920195895276243
$ws.Range("Q2").Value = 0.1316772311168889
$ws.Range("R2").Value = 1.185095080052
$ws.Range("S2").Value = 0.005301691708300269
$ws.Range("T2").Value = 0.005301691708300268

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il23a"
$ws.Range("C3").Value = "Il23r"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5663796666666666
$ws.Range("H3").Value = 1.699139
$ws.Range("I3").Value = 0.005761481588340165
$ws.Range("J3").Value = 0.005761481588340165
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02016266666666667
$ws.Range("N3").Value = 0.060488
$ws.Range("O3").Value = 0.07980410472375705
$ws.Range("P3").Value = 0.07980410472375705
$ws.Range("Q3").Value = 0.01141972442577778
$ws.Range("R3").Value = 0.102777519832
$ws.Range("S3").Value = 0.0004597898800398967
$ws.Range("T3").Value = 0.0004597898800398966

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Il23a"
$ws.Range("C4").Value = "Il23r"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2158466666666667
$ws.Range("H4").Value = 0.64754
$ws.Range("I4").Value = 0.00219569428264185
$ws.Range("J4").Value = 0.00219569428264185
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2324893333333333
$ws.Range("N4").Value = 0.697468
$ws.Range("O4").Value = 0.920195895276243
$ws.Range("P4").Value = 0.920195895276243
$ws.Range("Q4").Value = 0.05018204763555555
$ws.Range("R4").Value = 0.45163842872
$ws.Range("S4").Value = 0.002020468866168545
$ws.Range("T4").Value = 0.002020468866168545

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il23a"
$ws.Range("C5").Value = "Il23r"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2158466666666667
$ws.Range("H5").Value = 0.64754
$ws.Range("I5").Value = 0.00219569428264185
$ws.Range("J5").Value = 0.00219569428264185
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02016266666666667
$ws.Range("N5").Value = 0.060488
$ws.Range("O5").Value = 0.07980410472375705
$ws.Range("P5").Value = 0.07980410472375705
$ws.Range("Q5").Value = 0.004352044391111111
$ws.Range("R5").Value = 0.03916839952
$ws.Range("S5").Value = 0.0001752254164733048
$ws.Range("T5").Value = 0.0001752254164733048

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Il23a"
$ws.Range("C6").Value = "Il23r"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.910797666666667
$ws.Range("H6").Value = 8.732393
$ws.Range("I6").Value = 0.0296100092409453
$ws.Range("J6").Value = 0.02961000924094529
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2324893333333333
$ws.Range("N6").Value = 0.697468
$ws.Range("O6").Value = 0.920195895276243
$ws.Range("P6").Value = 0.920195895276243
$ws.Range("Q6").Value = 0.6767294089915555
$ws.Range("R6").Value = 6.090564680923999
$ws.Range("S6").Value = 0.02724700896260949
$ws.Range("T6").Value = 0.02724700896260948

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Il23a"
$ws.Range("C7").Value = "Il23r"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.910797666666667
$ws.Range("H7").Value = 8.732393
$ws.Range("I7").Value = 0.0296100092409453
$ws.Range("J7").Value = 0.02961000924094529
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02016266666666667
$ws.Range("N7").Value = 0.060488
$ws.Range("O7").Value = 0.07980410472375705
$ws.Range("P7").Value = 0.07980410472375705
$ws.Range("Q7").Value = 0.05868944308711111
$ws.Range("R7").Value = 0.528204987784
$ws.Range("S7").Value = 0.002363000278335813
$ws.Range("T7").Value = 0.002363000278335812

$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "Il23a"
$ws.Range("C8").Value = "Il23r"
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 94.008606
$ws.Range("H8").Value = 282.025818
$ws.Range("I8").Value = 0.9562999600642293
$ws.Range("J8").Value = 0.9562999600642292
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2324893333333333
$ws.Range("N8").Value = 0.697468
$ws.Range("O8").Value = 0.920195895276243
$ws.Range("P8").Value = 0.920195895276243
$ws.Range("Q8").Value = 21.855998136536
$ws.Range("R8").Value = 196.703983228824
$ws.Range("S8").Value = 0.8799832979039389
$ws.Range("T8").Value = 0.8799832979039388

$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "Il23a"
$ws.Range("C9").Value = "Il23r"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 94.008606
$ws.Range("H9").Value = 282.025818
$ws.Range("I9").Value = 0.9562999600642293
$ws.Range("J9").Value = 0.9562999600642292
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.02016266666666667
$ws.Range("N9").Value = 0.060488
$ws.Range("O9").Value = 0.07980410472375705
$ws.Range("P9").Value = 0.07980410472375705
$ws.Range("Q9").Value = 1.895464186576
$ws.Range("R9").Value = 17.059177679184
$ws.Range("S9").Value = 0.07631666216029045
$ws.Range("T9").Value = 0.07631666216029044

$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Il23a"
$ws.Range("C10").Value = "Il23r"
$ws.Range("D10").Value = "Neutrophils"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.6028873333333333
$ws.Range("H10").Value = 1.808662
$ws.Range("I10").Value = 0.00613285482384343
$ws.Range("J10").Value = 0.006132854823843428
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2324893333333333
$ws.Range("N10").Value = 0.697468
$ws.Range("O10").Value = 0.920195895276243
$ws.Range("P10").Value = 0.920195895276243
$ws.Range("Q10").Value = 0.1401648742017778
$ws.Range("R10").Value = 1.261483867816
$ws.Range("S10").Value = 0.005643427835225831
$ws.Range("T10").Value = 0.005643427835225829

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Il23a"
$ws.Range("C11").Value = "Il23r"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.6028873333333333
$ws.Range("H11").Value = 1.808662
$ws.Range("I11").Value = 0.00613285482384343
$ws.Range("J11").Value = 0.006132854823843428
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.02016266666666667
$ws.Range("N11").Value = 0.060488
$ws.Range("O11").Value = 0.07980410472375705
$ws.Range("P11").Value = 0.07980410472375705
$ws.Range("Q11").Value = 0.01215581633955556
$ws.Range("R11").Value = 0.109402347056
$ws.Range("S11").Value = 0.0004894269886175997
$ws.Range("T11").Value = 0.0004894269886175996

$ws.Rows.Item(13).Delete()
$ws.Rows.Item(12).Delete()
